# Insert a new column A ("k") in front of the existing i/j/d/travel time
# table, shifting the current data right by one column, then fill the new
# column with the value 1 for every data row (read from the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:D -> B:E by inserting a new blank column at A.
$ws.Range("A1:A13").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "k"

# New column values (constant 1 for every one of the 12 data rows).
$kValues = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $kValues[$i]
}

# Update the active selection to match the post-edit workbook state.
$ws.Range("E14").Select() | Out-Null
